# feat: add 2022-Q3 data
#
# - Insert a new worksheet "2022-Q3" positioned right after "总计" and
#   before "2022-Q1".
# - Populate it with the same layout/headers as the other quarterly
#   sheets (基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# - Update the "总计" summary sheet: insert a new leading data row for
#   2022-Q3 and shift the existing rows down.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    # Force the cell to be stored as text (matches cells like "004317" /
    # "1.79" / "0.0541" elsewhere in the workbook), then drop the
    # quote-prefix/number-format style back off so the cell ends up with
    # the default (no) style, same as its siblings.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet, inserted before "2022-Q1" (i.e.
#    right after "总计"), cloning the layout/styles of an existing
#    per-quarter sheet.
# ---------------------------------------------------------------------
$q1Sheet = $wb.Worksheets.Item("2022-Q1")
$q3Sheet = $wb.Worksheets.Add($q1Sheet)
$q3Sheet.Name = "2022-Q3"

$wb.Worksheets.Item("2021-Q4").Range("A1:H3").Copy($q3Sheet.Range("A1"))
$q3Sheet.Range("A1").ClearContents()

# Row 2 data.
$q3Sheet.Range("A2").Value = 0
Set-TextValue $q3Sheet.Range("B2") "004317"
Set-TextValue $q3Sheet.Range("C2") "前海开源沪港深裕鑫灵活配置混合C"
Set-TextValue $q3Sheet.Range("D2") "1.79"
Set-TextValue $q3Sheet.Range("E2") "70.17"
Set-TextValue $q3Sheet.Range("F2") "3.02"
Set-TextValue $q3Sheet.Range("G2") "0.0541"
$q3Sheet.Range("H2").Value = 3

# Row 3 data.
$q3Sheet.Range("A3").Value = 1
Set-TextValue $q3Sheet.Range("B3") "004316"
Set-TextValue $q3Sheet.Range("C3") "前海开源沪港深裕鑫灵活配置混合A"
Set-TextValue $q3Sheet.Range("D3") "1.77"
Set-TextValue $q3Sheet.Range("E3") "70.17"
Set-TextValue $q3Sheet.Range("F3") "3.02"
Set-TextValue $q3Sheet.Range("G3") "0.0535"
$q3Sheet.Range("H3").Value = 3

# ---------------------------------------------------------------------
# 2. Update "总计": shift the existing three rows down by one and add
#    the new 2022-Q3 row on top.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Range("A4:D4").Copy($total.Range("A5"))
$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q3"
$total.Range("C5").Value = 2
$total.Range("D5").Value = 0.04

$total.Range("B4").Value = "2021-Q4"
$total.Range("C4").Value = 2
$total.Range("D4").Value = 0.14

$total.Range("B3").Value = "2022-Q1"
$total.Range("C3").Value = 1
$total.Range("D3").Value = 0.11

$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.11
